$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 270; existing rows 270:395 shift down to 271:396
$ws.Rows(270).Insert()

# Populate the newly inserted row 270 with its data
$ws.Range("A270").Value = 9
$ws.Range("B270").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C270").Value = "Metropolitana"
$ws.Range("D270").Value = 44992
$ws.Range("E270").Value = 13
$ws.Range("F270").Value = 300000001
$ws.Range("G270").Value = "Rabanito"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 7000
$ws.Range("K270").Value = 3000
$ws.Range("L270").Value = 3000
$ws.Range("M270").Value = 3000
$ws.Range("N270").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O270").Value = "Provincia de Chacabuco"
$ws.Range("P270").Value = 30
$ws.Range("Q270").Value = 100
$ws.Range("R270").Value = "Hortaliza"
